$wb = $excel.ActiveWorkbook

# --- 1. Remove the "Parameters" worksheet entirely ---
$paramSheet = $wb.Worksheets.Item("Parameters")
$paramSheet.Delete()

# --- 2. Work on the remaining sheet ("Лист1") ---
$ws = $wb.Worksheets.Item(1)

# D4 currently holds "484" as text; make it a genuine number like the rest
# of the numeric column.
$ws.Range("D4").Value = 484

# --- 3. Extend the table with yearly columns E:H (2020-2023), copying the
#         number formatting from column D (the 2019 column) row by row ---
$ws.Range("D3:D6").Copy()
$ws.Range("E3:H6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 3: year headers
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# Row 4: number of local governments (constant across years)
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# Row 5: proportion of local governments (%)
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# Row 6: number of local governments with DRR strategies
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169

# --- 4. Restore the cursor position left by the original author ---
$ws.Range("D9").Select()
